$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the cells that hold numeric-looking or percentage-looking text values to
# the Text number format first, so Excel keeps the assigned value as a literal
# string (matching the inlineStr cells in the source workbook) instead of
# re-interpreting it as a number/percentage.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "298.16"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "-1.03%"
$c = $ws.Range("G2")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "31.77"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "0.95%"
$c = $ws.Range("G3")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "5.083"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "-0.90%"
$c = $ws.Range("G4")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.08114"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "9.67%"
$c = $ws.Range("G5")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "2.591"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "7.93%"
$c = $ws.Range("G6")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "7.756"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "-2.17%"
$c = $ws.Range("G7")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "3.830"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "1.67%"
$c = $ws.Range("G8")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "0.17%"
$c = $ws.Range("G9")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.1758"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "1.71%"
$c = $ws.Range("G10")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07483"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "-1.73%"
$c = $ws.Range("G11")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.08923"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "10.66%"
$c = $ws.Range("G12")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.03036"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "-0.02%"
$c = $ws.Range("G13")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "0.87%"
$c = $ws.Range("G14")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.001495"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "-0.52%"
$c = $ws.Range("G15")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.006048"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "-0.52%"
$c = $ws.Range("G16")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.566"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "3.06%"
$c = $ws.Range("G17")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.253"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "1.05%"
$c = $ws.Range("G18")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "-0.62%"
$c = $ws.Range("G19")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.1316"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "-0.67%"
$c = $ws.Range("G20")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "3.971"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "-14.60%"
$c = $ws.Range("G21")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.1697"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "7.15%"
$c = $ws.Range("G22")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.04591"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "-1.45%"
$c = $ws.Range("G23")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.001243"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "1.46%"
$c = $ws.Range("G24")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.004466"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "-0.60%"
$c = $ws.Range("G25")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.0001198"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "-7.83%"
$c = $ws.Range("G26")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.0003407"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "82.01%"
$c = $ws.Range("G27")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("G28")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("G29")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("G30")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("G31")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("G32")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("G33")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("G34")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("G35")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("G36")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("G37")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("G38")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01775"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "3.23%"
$c = $ws.Range("G39")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.04514"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "-0.34%"
$c = $ws.Range("G40")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.006895"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "-2.51%"
$c = $ws.Range("G41")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1354"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "0.58%"
$c = $ws.Range("G42")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.002206"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "-1.04%"
$c = $ws.Range("G43")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.009854"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "-9.91%"
$c = $ws.Range("G44")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.00006109"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "-2.58%"
$c = $ws.Range("G45")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.00000000749"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "-0.19%"
$c = $ws.Range("G46")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.008734"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "24.83%"
$c = $ws.Range("G47")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "-55.60%"
$c = $ws.Range("G48")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.00002096"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "-0.19%"
$c = $ws.Range("G49")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0001996"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "-0.12%"
$c = $ws.Range("G50")
$c.NumberFormat = "@"
$c.Value = "18"

$c = $ws.Range("G51")
$c.NumberFormat = "@"
$c.Value = "18"
